$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.294.35'
$ws.Range("E2").Value = '  -1.40%  '

$ws.Range("D3").Value = '2.995.84'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.62'
$ws.Range("E5").Value = '  +1.25%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.89'
$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.431'
$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("E9").Value = '  -2.38%  '

$ws.Range("E10").Value = '  -0.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.367'
$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("D12").Value = '3.507.82'
$ws.Range("E12").Value = '  +0.37%  '

$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.51'
$ws.Range("E14").Value = '  -2.17%  '

$ws.Range("E15").Value = '  +1.97%  '

$ws.Range("D16").Value = '56.264.09'
$ws.Range("E16").Value = '  -1.51%  '

$ws.Range("D17").Value = '2.989.43'
$ws.Range("E17").Value = '  +0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.97'
$ws.Range("E18").Value = '  -1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.91'
$ws.Range("E19").Value = '  +2.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.03'
$ws.Range("E20").Value = '  +1.99%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.89'
$ws.Range("E21").Value = '  +3.52%  '

$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.497'
$ws.Range("E23").Value = '  +0.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.97'
$ws.Range("E24").Value = '  +3.17%  '

$ws.Range("D25").Value = '3.121.11'
$ws.Range("E25").Value = '  +0.44%  '

$ws.Range("E26").Value = '  +1.45%  '

$ws.Range("E27").Value = '  -0.73%  '

$ws.Range("D28").Value = '0.0₃0936'
$ws.Range("E28").Value = '  +4.94%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.33'
$ws.Range("E29").Value = '  -4.31%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.89'
$ws.Range("E30").Value = '  -3.09%  '

$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("E33").Value = '  -0.57%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.07'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.46'
$ws.Range("E35").Value = '  -2.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.84'
$ws.Range("E36").Value = '  +0.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.25'
$ws.Range("E37").Value = '  +7.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.23'
$ws.Range("E38").Value = '  -0.94%  '

$ws.Range("E39").Value = '  -0.42%  '

$ws.Range("D40").Value = '3.032.43'
$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("E43").Value = '  +1.40%  '

$ws.Range("E44").Value = '  +1.02%  '

$ws.Range("D45").Value = '2.185.49'
$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("E46").Value = '  -2.73%  '

$ws.Range("E47").Value = '  -1.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.928'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0237'
$ws.Range("E49").Value = '  +1.19%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.39'
$ws.Range("E50").Value = '  +1.12%  '

$ws.Range("E51").Value = '  -1.94%  '
